$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume (E) figures for the latest data pull.
# Rows 48 and 49 also swap their Coin/Link (B/C) content (ranking order changed).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.440.39"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.724.43"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.37"
$ws.Range("E5").Value = "  -1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4881"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("E8").Value = "  -2.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06211"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.713.07"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06991"
$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.49"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.549"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5989"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.41"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.445.54"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007302"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.942.97"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.469"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.534"
$ws.Range("E23").Value = "  -2.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.136"
$ws.Range("E24").Value = "  -2.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.14"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.28"
$ws.Range("E26").Value = "  -0.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "107.12"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.730"
$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.960"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07979"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.677"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04508"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.602"
$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6304"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9319"
$ws.Range("E37").Value = "  +4.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.967"
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.390"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.84"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.342"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3849"
$ws.Range("E44").Value = "  -1.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.875"
$ws.Range("E45").Value = "  -2.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1171"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05360"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.29"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.729"
$ws.Range("E49").Value = "  -2.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.235"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.98"
$ws.Range("E51").Value = "  -0.85%  "
